# Added handling of common packages.
# Reorders the data rows of the "classFields" worksheet so that fields
# belonging to the same class are grouped/ordered consistently (e.g. enum
# constants together, CORS-filter fields together, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# New row order (Class Name, Field Name, Field Modifier, Field Type),
# rows 2-19 -- header row (1) is left untouched.
$rows = @(
    @{Row=2;  A="org.andante.config.security.role.KeycloakRole"; B='$VALUES'; C="private"; D="org.andante.config.security.role.KeycloakRole[]"},
    @{Row=3;  A="org.andante.config.security.role.KeycloakRole"; B="BLOGGER"; C="public"; D="org.andante.config.security.role.KeycloakRole"},
    @{Row=4;  A="org.andante.config.security.role.KeycloakRole"; B="ADMIN"; C="public"; D="org.andante.config.security.role.KeycloakRole"},
    @{Row=5;  A="org.andante.config.security.role.KeycloakRole"; B="name"; C="private"; D="java.lang.String"},
    @{Row=6;  A="org.andante.config.gateway.GatewayConfiguration"; B="filterFactory"; C="private"; D="org.springframework.cloud.gateway.filter.factory.TokenRelayGatewayFilterFactory"},
    @{Row=7;  A="org.andante.config.security.filter.CrossOriginRequestSharingFilter"; B="allowedHeaders"; C="private"; D="java.lang.String"},
    @{Row=8;  A="org.andante.config.security.filter.CrossOriginRequestSharingFilter"; B="allowedOrigins"; C="private"; D="java.lang.String"},
    @{Row=9;  A="org.andante.config.security.filter.CrossOriginRequestSharingFilter"; B="allowedMethods"; C="private"; D="java.lang.String"},
    @{Row=10; A="org.andante.config.security.filter.CrossOriginRequestSharingFilter"; B="exposedHeaders"; C="private"; D="java.lang.String"},
    @{Row=11; A="org.andante.config.security.converter.KeycloakRealmRoleConverter"; B="ROLES"; C="private"; D="java.lang.String"},
    @{Row=12; A="org.andante.config.security.converter.KeycloakRealmRoleConverter"; B="REALM_ACCESS"; C="private"; D="java.lang.String"},
    @{Row=13; A="org.andante.config.security.SecurityConfiguration"; B="exposedHeaders"; C="private"; D="java.util.List"},
    @{Row=14; A="org.andante.config.security.SecurityConfiguration"; B="keycloakRealmRoleConverter"; C="private"; D="org.andante.config.security.converter.KeycloakRealmRoleConverter"},
    @{Row=15; A="org.andante.config.security.SecurityConfiguration"; B="allowedHeaders"; C="private"; D="java.util.List"},
    @{Row=16; A="org.andante.config.security.SecurityConfiguration"; B="disabledSecurityEndpoints"; C="private"; D="java.util.List"},
    @{Row=17; A="org.andante.config.security.SecurityConfiguration"; B="jwkSetUri"; C="private"; D="java.lang.String"},
    @{Row=18; A="org.andante.config.security.SecurityConfiguration"; B="allowedMethods"; C="private"; D="java.util.List"},
    @{Row=19; A="org.andante.config.security.SecurityConfiguration"; B="allowedOrigins"; C="private"; D="java.util.List"}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
